$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 9 new rows at row 39 (pushes the existing 2019-11-29.. data down by 9 rows,
# from rows 39:111 to rows 48:120), to make room for the previously-missing
# trading days 2019-11-18 .. 2019-11-28.
$ws.Range("A39:A47").EntireRow.Insert()

# Columns B, C and D hold text that looks numeric/date-like ("2019-11-18", "0213").
# Force a text format first so Excel doesn't auto-convert them into date serials
# or numbers (which would drop the leading zero on "0213"), then restore the
# default "Normal" style afterwards so no stray formatting is left behind.
$ws.Range("B39:D47").NumberFormat = "@"

$ws.Range("A39").Value = 1574035200
$ws.Range("B39").Value = "2019-11-18"
$ws.Range("C39").Value = "0213"
$ws.Range("D39").Value = "MTAG"
$ws.Range("E39").Value = 0.6
$ws.Range("F39").Value = 0.605
$ws.Range("G39").Value = 0.58
$ws.Range("H39").Value = 0.59
$ws.Range("I39").Value = 21923800

$ws.Range("A40").Value = 1574121600
$ws.Range("B40").Value = "2019-11-19"
$ws.Range("C40").Value = "0213"
$ws.Range("D40").Value = "MTAG"
$ws.Range("E40").Value = 0.585
$ws.Range("F40").Value = 0.59
$ws.Range("G40").Value = 0.5600000000000001
$ws.Range("H40").Value = 0.57
$ws.Range("I40").Value = 23531000

$ws.Range("A41").Value = 1574208000
$ws.Range("B41").Value = "2019-11-20"
$ws.Range("C41").Value = "0213"
$ws.Range("D41").Value = "MTAG"
$ws.Range("E41").Value = 0.5649999999999999
$ws.Range("F41").Value = 0.6
$ws.Range("G41").Value = 0.555
$ws.Range("H41").Value = 0.585
$ws.Range("I41").Value = 51259400

$ws.Range("A42").Value = 1574294400
$ws.Range("B42").Value = "2019-11-21"
$ws.Range("C42").Value = "0213"
$ws.Range("D42").Value = "MTAG"
$ws.Range("E42").Value = 0.585
$ws.Range("F42").Value = 0.59
$ws.Range("G42").Value = 0.5649999999999999
$ws.Range("H42").Value = 0.57
$ws.Range("I42").Value = 16442100

$ws.Range("A43").Value = 1574380800
$ws.Range("B43").Value = "2019-11-22"
$ws.Range("C43").Value = "0213"
$ws.Range("D43").Value = "MTAG"
$ws.Range("E43").Value = 0.57
$ws.Range("F43").Value = 0.585
$ws.Range("G43").Value = 0.5649999999999999
$ws.Range("H43").Value = 0.575
$ws.Range("I43").Value = 13173400

$ws.Range("A44").Value = 1574640000
$ws.Range("B44").Value = "2019-11-25"
$ws.Range("C44").Value = "0213"
$ws.Range("D44").Value = "MTAG"
$ws.Range("E44").Value = 0.58
$ws.Range("F44").Value = 0.58
$ws.Range("G44").Value = 0.54
$ws.Range("H44").Value = 0.54
$ws.Range("I44").Value = 19477400

$ws.Range("A45").Value = 1574726400
$ws.Range("B45").Value = "2019-11-26"
$ws.Range("C45").Value = "0213"
$ws.Range("D45").Value = "MTAG"
$ws.Range("E45").Value = 0.54
$ws.Range("F45").Value = 0.5600000000000001
$ws.Range("G45").Value = 0.525
$ws.Range("H45").Value = 0.525
$ws.Range("I45").Value = 14542200

$ws.Range("A46").Value = 1574812800
$ws.Range("B46").Value = "2019-11-27"
$ws.Range("C46").Value = "0213"
$ws.Range("D46").Value = "MTAG"
$ws.Range("E46").Value = 0.525
$ws.Range("F46").Value = 0.54
$ws.Range("G46").Value = 0.52
$ws.Range("H46").Value = 0.525
$ws.Range("I46").Value = 8898800

$ws.Range("A47").Value = 1574899200
$ws.Range("B47").Value = "2019-11-28"
$ws.Range("C47").Value = "0213"
$ws.Range("D47").Value = "MTAG"
$ws.Range("E47").Value = 0.525
$ws.Range("F47").Value = 0.53
$ws.Range("G47").Value = 0.505
$ws.Range("H47").Value = 0.515
$ws.Range("I47").Value = 12177600

# Restore default styling on the text columns now that the values are locked in as text.
$ws.Range("B39:D47").Style = "Normal"
